# Applies the cryptos list price/volume refresh described in the commit
# "Updated cryptos list on Tue Nov 28 11:49:46 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing formatting such as trailing zeros e.g. "1.00" -> 1, "0.610" -> 0.61).
# Temporarily force the Text format so the literal string is kept, then restore
# the default "Normal" style so no stray formatting is left behind.
$textCells = @("D5", "D6", "D8", "D10", "D13", "D15", "D16", "D24", "D26", "D27", "D29", "D35", "D36", "D37", "D38", "D42", "D46")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.198.80'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '2.025.63'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '228.68'
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("D6").Value = '0.610'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '56.12'
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '0.0783'
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").Value = '2.325.96'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '14.31'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("D15").Value = '0.739'
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("D16").Value = '5.19'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = '2.022.51'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '37.184.12'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '2.44'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").Value = '163.45'
$ws.Range("E26").Value = '  -2.17%  '
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -4.10%  '
$ws.Range("E28").Value = '  +1.97%  '
$ws.Range("D29").Value = '18.75'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.93'
$ws.Range("E35").Value = '  +5.39%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '2.36'
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '3.20'
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("E39").Value = '  +2.10%  '
$ws.Range("D40").Value = '1.473.61'
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("E41").Value = '  -2.13%  '
$ws.Range("D42").Value = '94.29'
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("E44").Value = '  -1.85%  '
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").Value = '4.18'
$ws.Range("E46").Value = '  +13.51%  '
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").Value = '2.213.40'
$ws.Range("E51").Value = '  +0.03%  '

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
